$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "tracking_number"
$ws.Range("B1").Value = "carrier_slug"
$ws.Range("C1").Value = "status_tag"
$ws.Range("D1").Value = "title"
$ws.Range("E1").Value = "order_id"
$ws.Range("F1").Value = "last_checkpoint_time"
$ws.Range("G1").Value = "last_checkpoint_location"
$ws.Range("H1").Value = "updated_at"

# Row 2
$ws.Range("A2").Value = "TEST_GDPR"
$ws.Range("B2").Value = "dbschenker-se"
$ws.Range("C2").Value = "Delivered"
$ws.Range("D2").Value = "SHIPMENT_TITLE"
$ws.Range("H2").Value = "2026-02-08T12:16:34+00:00"

# Row 3
$ws.Range("A3").Value = "TEST_TRACKING_DECEMBER"
$ws.Range("B3").Value = "dhl"
$ws.Range("C3").Value = "Delivered"
$ws.Range("D3").Value = "TEST_Tracking_December"
$ws.Range("H3").Value = "2026-02-08T11:24:55+00:00"

# Row 4
$ws.Range("A4").Value = "TEST_TRACKING"
$ws.Range("B4").Value = "kn"
$ws.Range("C4").Value = "Delivered"
$ws.Range("D4").Value = "ZFRE"
$ws.Range("H4").Value = "2026-02-08T11:24:18+00:00"

# Row 5
$ws.Range("A5").Value = "ITD-0-12345678"
$ws.Range("B5").Value = "testing-courier"
$ws.Range("C5").Value = "Delivered"
$ws.Range("D5").Value = "ITD-0-12345678"
$ws.Range("H5").Value = "2026-02-08T11:07:24+00:00"
